# Add data for 2022-10-21
# Updates the 2022 (column I) violent-crime counts across the citywide
# totals sheet, the by-neighborhood summary sheet, and each individual
# neighborhood sheet to reflect the newly ingested day of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5896
$ws.Range("I3").Value = 6125
$ws.Range("I4").Value = 1404
$ws.Range("I5").Value = 564
$ws.Range("I6").Value = 6920
$ws.Range("I7").Value = 20909

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 59
$ws.Range("I7").Value = 241

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 200
$ws.Range("I3").Value = 222
$ws.Range("I6").Value = 193
$ws.Range("I7").Value = 668

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 298
$ws.Range("I7").Value = 810

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 208

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 66
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 162
$ws.Range("I6").Value = 150
$ws.Range("I7").Value = 658
$ws.Range("I8").Value = 1255
$ws.Range("I11").Value = 312
$ws.Range("I12").Value = 47
$ws.Range("I15").Value = 234
$ws.Range("I18").Value = 154
$ws.Range("I19").Value = 577
$ws.Range("I20").Value = 517
$ws.Range("I21").Value = 99
$ws.Range("I23").Value = 210
$ws.Range("I25").Value = 110
$ws.Range("I29").Value = 1303
$ws.Range("I31").Value = 208
$ws.Range("I33").Value = 949
$ws.Range("I36").Value = 279
$ws.Range("I37").Value = 668
$ws.Range("I40").Value = 36
$ws.Range("I41").Value = 87
$ws.Range("I42").Value = 711
$ws.Range("I44").Value = 154
$ws.Range("I48").Value = 284
$ws.Range("I53").Value = 217
$ws.Range("I54").Value = 431
$ws.Range("I55").Value = 230
$ws.Range("I63").Value = 74
$ws.Range("I67").Value = 810
$ws.Range("I72").Value = 83
$ws.Range("I73").Value = 191
$ws.Range("I76").Value = 302
$ws.Range("I77").Value = 133
$ws.Range("I79").Value = 589
$ws.Range("I83").Value = 451
$ws.Range("I84").Value = 182
$ws.Range("I85").Value = 950
$ws.Range("I86").Value = 130
$ws.Range("I87").Value = 49
$ws.Range("I88").Value = 189
$ws.Range("I89").Value = 241
$ws.Range("I90").Value = 254
$ws.Range("I97").Value = 175
$ws.Range("I101").Value = 20909

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 155
$ws.Range("I5").Value = 20
$ws.Range("I6").Value = 97
$ws.Range("I7").Value = 451

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 360
$ws.Range("I6").Value = 299
$ws.Range("I7").Value = 949

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 96
$ws.Range("I7").Value = 431

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 381
$ws.Range("I3").Value = 450
$ws.Range("I4").Value = 67
$ws.Range("I7").Value = 1303

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 197
$ws.Range("I3").Value = 177
$ws.Range("I6").Value = 168
$ws.Range("I7").Value = 577

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 48
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I3").Value = 55
$ws.Range("I7").Value = 284

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 57
$ws.Range("I6").Value = 141
$ws.Range("I7").Value = 302

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 269
$ws.Range("I3").Value = 365
$ws.Range("I6").Value = 238
$ws.Range("I7").Value = 950

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 150

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 179
$ws.Range("I3").Value = 234
$ws.Range("I5").Value = 24
$ws.Range("I6").Value = 225
$ws.Range("I7").Value = 711

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 73
$ws.Range("I7").Value = 230

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 73
$ws.Range("I7").Value = 210

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 171
$ws.Range("I6").Value = 173
$ws.Range("I7").Value = 589

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 146
$ws.Range("I3").Value = 146
$ws.Range("I6").Value = 179
$ws.Range("I7").Value = 517

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 82
$ws.Range("I3").Value = 90
$ws.Range("I7").Value = 279

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 53
$ws.Range("I7").Value = 234

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 130
$ws.Range("I7").Value = 312

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 56
$ws.Range("I7").Value = 162

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 175

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I3").Value = 66
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 386
$ws.Range("I7").Value = 1255

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I3").Value = 9
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 83
$ws.Range("I7").Value = 254

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 217

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I2").Value = 17
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 83

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I3").Value = 46
$ws.Range("I7").Value = 133

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 205
$ws.Range("I6").Value = 174
$ws.Range("I7").Value = 658

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 47

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 49
